# Updated cryptos list refresh.
# Applies new Price (column D) and Volume(1h) (column E) values, and
# swaps the FraxShare / RenderToken rows (43 <-> 44) with their own new
# Price / Volume(1h) figures.
#
# Column D holds price strings that can look like plain numbers (e.g.
# "23.95"); assigning them straight to .Value lets Excel auto-convert
# them to numeric cells, which changes the stored cell type. Prefixing
# with a literal apostrophe forces text entry (same as typing '23.95
# into Excel), then resetting .Style to "Normal" drops the resulting
# quote-prefix cell style so the cell ends up as plain, unstyled text -
# matching the original file's cells (no style index, inline text).

function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "28.209.00"
$ws.Range("E2").Value = "  +2.36%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.588.60"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.86%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "213.12"
$ws.Range("E5").Value = "  +0.69%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.26%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.96%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "23.95"
$ws.Range("E8").Value = "  +6.18%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.250"
$ws.Range("E9").Value = "  -0.59%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.31%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +2.52%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") "1.817.00"
$ws.Range("E12").Value = "  +1.08%  "

# Row 13 - WrappedEther
Set-TextValue $ws.Range("D13") "1.590.96"
$ws.Range("E13").Value = "  +1.31%  "

# Row 14 - Polygon
Set-TextValue $ws.Range("D14") "0.530"
$ws.Range("E14").Value = "  +1.08%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  -1.12%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "28.263.18"
$ws.Range("E16").Value = "  +2.68%  "

# Row 17 - Litecoin
$ws.Range("E17").Value = "  +1.88%  "

# Row 18 - BitcoinCash
Set-TextValue $ws.Range("D18") "227.06"
$ws.Range("E18").Value = "  +0.61%  "

# Row 19 - ShibaInu (subscript-3 digit in the price string)
$d19value = "{0}{1}{2}" -f "0.0", [char]8323, "0708"
Set-TextValue $ws.Range("D19") $d19value
$ws.Range("E19").Value = "  +0.11%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -1.29%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.86%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -1.99%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  -1.26%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.35%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "152.03"
$ws.Range("E25").Value = "  +0.83%  "

# Row 26 - EthereumClassic
$ws.Range("E26").Value = "  -0.18%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.48%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  -1.39%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.93%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.47%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.16%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.68%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -1.42%  "

# Row 34 - Maker
Set-TextValue $ws.Range("D34") "1.397.18"
$ws.Range("E34").Value = "  -4.12%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -2.81%  "

# Row 36 - TrustWalletToken
Set-TextValue $ws.Range("D36") "1.02"
$ws.Range("E36").Value = "  -8.28%  "

# Row 37 - HuobiToken
Set-TextValue $ws.Range("D37") "2.36"
$ws.Range("E37").Value = "  +1.38%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.81%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  +8.32%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "  -0.23%  "

# Row 41 - ARBITRUM
Set-TextValue $ws.Range("D41") "0.813"
$ws.Range("E41").Value = "  -0.82%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.77%  "

# Row 43 / 44 - FraxShare and RenderToken swap ranking positions
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D43") "5.59"
$ws.Range("E43").Value = "  -3.88%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D44") "1.87"
$ws.Range("E44").Value = "  +4.04%  "

# Row 45 - WEMIXToken
Set-TextValue $ws.Range("D45") "0.981"
$ws.Range("E45").Value = "  +1.17%  "

# Row 46 - Aave
Set-TextValue $ws.Range("D46") "64.29"
$ws.Range("E46").Value = "  -2.00%  "

# Row 47 - RocketPoolETH
Set-TextValue $ws.Range("D47") "1.726.40"
$ws.Range("E47").Value = "  +0.87%  "

# Row 48 - Quant
Set-TextValue $ws.Range("D48") "87.35"
$ws.Range("E48").Value = "  +0.90%  "

# Row 49 - mCoin
$ws.Range("E49").Value = "  +1.63%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  +6.58%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -0.70%  "
